# Update automatico via Actualizar 05-23-2020 07-32-32
# Adds the newest COVID-19 condition-of-patients record (2020-05-22) as a
# new row at the bottom of the "Condicion_Pacientes" table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Grow the table by one row - this also keeps the table ref / autoFilter
# ref and the sheet dimension in sync with the new row.
$newRow = $lo.ListRows.Add()

# Match the formatting of the row above it (date format on column A,
# centered alignment on the rest) without minting new style entries.
$ws.Range("A70:F70").Copy()
$ws.Range("A71:F71").PasteSpecial(-4122)

# Fill in the new record (row 71): 2020-05-22
$ws.Range("A71").Value = 43973
$ws.Range("B71").Value = 592
$ws.Range("C71").Value = 273
$ws.Range("D71").Value = 347
$ws.Range("E71").Value = 14
$ws.Range("F71").Value = 20

# Move the active selection to the newly added cell, matching the
# author's final cursor position.
$null = $ws.Range("C71").Select()
